$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the volatile RANDBETWEEN-based access-code formulas in column A
# with their static (already-computed) values.
$ws.Range("A3").Value  = "DERPEN002"
$ws.Range("A4").Value  = "RAYMEA005"
$ws.Range("A5").Value  = "KELBON006"
$ws.Range("A6").Value  = "INDLAM001"
$ws.Range("A7").Value  = "EDIROM005"
$ws.Range("A8").Value  = "JAVBLA001"
$ws.Range("A9").Value  = "QUIHOU009"
$ws.Range("A10").Value = "ALVTUC007"
$ws.Range("A11").Value = "PHIPAL006"
$ws.Range("A12").Value = "ALICHA008"
$ws.Range("A13").Value = "JAVBUR009"

# New ip_address column (D)
$ws.Range("D1").Value  = "ip_address"
$ws.Range("D2").Value  = "88.60.241.111"
$ws.Range("D3").Value  = "2.157.164.237"
$ws.Range("D4").Value  = "181.87.13.187"
$ws.Range("D5").Value  = "187.160.100.85"
$ws.Range("D6").Value  = "62.182.139.127"
$ws.Range("D7").Value  = "27.158.255.22"
$ws.Range("D8").Value  = "117.105.21.50"
$ws.Range("D9").Value  = "8.248.52.152"
$ws.Range("D10").Value = "7.190.48.73"
$ws.Range("D11").Value = "153.53.191.20"
$ws.Range("D12").Value = "24.29.158.165"
$ws.Range("D13").Value = "21.184.96.247"

# Mirror the saved selection/view state from the authored workbook
[void]$ws.Range("K10").Select()
